$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsParameters = $wb.Worksheets.Item("parameters")

$wsParameters.Range("B2").Value = 100
$wsParameters.Range("B3").Value = 0.05
$wsParameters.Range("B6").Value = 1000
$wsParameters.Range("B8").Value = "1 0.5 0.2 0.01"

$wsInstructions.Activate()
$wsInstructions.Range("A12").Select() | Out-Null
$wsParameters.Activate()
$wsParameters.Range("B8").Select() | Out-Null
